$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in the user-story texts
$ws.Range("A3").Value = "As a player I want my party to get stronger as I progress through the game"
$ws.Range("A18").Value = "As a player I want the soulbar mechanic to make myself and the enemies stronger as well as provide interesting situations"
$ws.Range("A22").Value = "As a player I want to choose my characters properties at the beginning of the game"

# Fill in the new "ID" numbers in column B (Burn, Freeze, Paralyzed, etc.)
$ws.Range("B2").Value = 1
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 3
$ws.Range("B15").Value = 4
$ws.Range("B25").Value = 5

# Reflect the editor's final selection state
$null = $ws.Range("B2").Select()
